# fix units for tax credits
#
# - "discount rate" header (col E, shared by both sheets) becomes "discount rate (-)"
# - investment_tax_credits "tax credit ($/output unit)" (col B) becomes
#   "tax credit (% inv. cost)"
# - investment_tax_credits "levelized tax credit ($/output unit)" (col I) becomes
#   "levelized tax credit (% inv. cost)"
# - production_tax_credits becomes the active sheet/tab (was investment_tax_credits)
# - both sheet views rezoom from 130% to 110% and move the selection

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("production_tax_credits")
$ws2 = $wb.Worksheets.Item("investment_tax_credits")

# Relabel the shared "discount rate" header on both sheets so the unit is explicit.
$ws1.Range("E1").Value = "discount rate (-)"
$ws2.Range("E1").Value = "discount rate (-)"

# Fix the investment tax credit sheet's unit labels: they're a % of investment
# cost, not a $/output-unit value.
$ws2.Range("B1").Value = "tax credit (% inv. cost)"
$ws2.Range("I1").Value = "levelized tax credit (% inv. cost)"

# investment_tax_credits: no longer the selected tab; zoom out a bit; move selection.
$ws2.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 110
$ws2.Range("E1").Select() | Out-Null

# production_tax_credits: becomes the selected tab; zoom out a bit; move selection.
$ws1.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 110
$ws1.Range("E2").Select() | Out-Null
